# Add data-dictionary entries for the "seedID" (recruitment) and the
# "sorted_by" (sp biomass) tables to the table_description sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- recruitment / seedID table -------------------------------------------------
$ws.Cells.Item(62, 2).Value = "seedID"

$ws.Cells.Item(63, 1).Value = "recruitment"
$ws.Cells.Item(63, 2).Value = "round"
$ws.Cells.Item(63, 3).Value = "Round of sampling; round 1-2 correspond to 2018, and round 3-4 to 2019"
$ws.Cells.Item(63, 5).Value = "defined"

$ws.Cells.Item(64, 2).Value = "presence"
$ws.Cells.Item(64, 3).Value = "Presence (1) or absence (0) of seedlings"
$ws.Cells.Item(64, 5).Value = "recorded"

$ws.Cells.Item(65, 2).Value = "x"
$ws.Cells.Item(65, 3).Value = "x coordinate in the plot"
$ws.Cells.Item(65, 5).Value = "recorded"

$ws.Cells.Item(66, 2).Value = "y"
$ws.Cells.Item(66, 3).Value = "y coordinate in the plot"
$ws.Cells.Item(66, 5).Value = "recorded"

# --- sp biomass / sorted_by -------------------------------------------------
$ws.Cells.Item(67, 2).Value = "sorted_by"
$ws.Cells.Item(67, 3).Value = "Person that sorted the biomass"
$ws.Cells.Item(67, 5).Value = "recorded"

# --- restore view state (top-left cell / selection / window geometry) -----
$null = $ws.Range("B68").Select()

$win = $excel.ActiveWindow
$win.ScrollRow = 40
